# EPICP-1: changed unit from mg/d to g/d for sodium and potassium intake
# in the DD_EPICP_INES data dictionary (Variables sheet).
#
# Row 25 (name "mna") held "natrium intake at baseline [mg/d]" and is
# renamed/re-unitized to "sodium intake at baseline [g/d]".
# Row 26 (name "mk") held "potassium intake at baseline [mg/d]" and is
# re-unitized to "potassium intake at baseline [g/d]".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Variables")
$ws.Activate()

# Update row 26 first, then row 25, so that the shared-string table ends up
# with the same index assignment as the target file.
$ws.Range("C26").Value2 = "potassium intake at baseline [g/d]"
$ws.Range("C25").Value2 = "sodium intake at baseline [g/d]"

# Leave the selection on C25, matching the saved worksheet view.
$ws.Range("C25").Select() | Out-Null
